$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update curated dimension/measure metadata (row 2)
$ws.Range("B2").Value = "iaest-measure:estado-civil"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:sexo"

# Update dim/medida labels (row 3) - columns B and D are now "medida"; column C is now "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"

# Update type labels (row 4) - columns B and D drop skos:Concept for xsd:int; column C becomes URI-Municipio
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("D4").Value = "xsd:int"

# Remove row 5 entirely (mapping-estado-civil.xlsx / mapping-sexo.xlsx no longer referenced)
$ws.Rows("5").Delete()
